# Update "每十万人口各级学校平均在校生数" data:
# Drop the oldest decade (2000-2009), keep 2010-2021 as-is, add new row for 2022,
# and refresh the 2020 row with its newly-computed (non-rounded) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2010年", 3955, 7448, 2230, 3504, 2189),
    @("2011年", 3779, 7403, 2554, 3495, 2253),
    @("2012年", 3535, 7196, 2736, 3411, 2335),
    @("2013年", 3279, 6913, 2876, 3227, 2418),
    @("2014年", 3222, 6946, 2977, 3100, 2488),
    @("2015年", 3152, 7086, 3118, 2965, 2524),
    @("2016年", 3150, 7211, 3211, 2887, 2530),
    @("2017年", 3213, 7300, 3327, 2861, 2576),
    @("2018年", 3347, 7438, 3350, 2828, 2658),
    @("2019年", 3459.3703507288, 7568.7166219954, 3378.2059367341, 2849.5681463114, 2856.8975476214),
    @("2020年", 3510, 7661, 3441, 2948, 3126),
    @("2021年", 3554, 7634, 3403, 2774, 3301),
    @("2022年", 3624.94442871301, 7597.38029166077, 3275.90867903157, 2895.15616593516, 3510.14222001982)
)

# The old sheet has 22 rows of data (rows 2-22); the new sheet only has 14
# (rows 2-14). Clear the old data range first, then delete the now-unused
# trailing rows so the sheet dimension shrinks to A1:F14.
$ws.Range("A2:F22").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove rows 15-22 entirely so the used range shrinks back down.
$ws.Range("A15:F22").EntireRow.Delete()
